$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '27.991.98'; E = '  +3.37%  '; DNumeric = $false },
    @{ Row = 3; D = '1.801.78'; E = '  +4.01%  '; DNumeric = $false },
    @{ Row = 4; D = '0.9986'; E = '  -0.22%  '; DNumeric = $true },
    @{ Row = 5; D = '315.35'; E = '  +1.52%  '; DNumeric = $true },
    @{ Row = 6; D = '0.9985'; E = '  -0.15%  '; DNumeric = $true },
    @{ Row = 7; D = '0.5424'; E = '  +11.23%  '; DNumeric = $true },
    @{ Row = 8; D = '0.3791'; E = '  +7.73%  '; DNumeric = $true },
    @{ Row = 9; D = '43.13'; E = '  -1.29%  '; DNumeric = $true },
    @{ Row = 10; D = '0.07563'; E = '  +3.78%  '; DNumeric = $true },
    @{ Row = 11; D = $null; E = '  +7.53%  '; DNumeric = $false },
    @{ Row = 12; D = '0.9983'; E = '  -0.24%  '; DNumeric = $true },
    @{ Row = 13; D = '21.11'; E = '  +5.37%  '; DNumeric = $true },
    @{ Row = 14; D = '6.228'; E = '  +5.61%  '; DNumeric = $true },
    @{ Row = 15; D = '1.797.24'; E = '  +3.56%  '; DNumeric = $false },
    @{ Row = 16; D = '7.141'; E = '  +3.34%  '; DNumeric = $true },
    @{ Row = 17; D = '91.75'; E = '  +4.87%  '; DNumeric = $true },
    @{ Row = 18; D = '0.00001080'; E = '  +3.75%  '; DNumeric = $true },
    @{ Row = 19; D = '0.06503'; E = '  +1.40%  '; DNumeric = $true },
    @{ Row = 20; D = '0.9980'; E = '  -0.19%  '; DNumeric = $true },
    @{ Row = 21; D = '17.14'; E = '  +3.01%  '; DNumeric = $true },
    @{ Row = 22; D = '5.974'; E = '  +4.61%  '; DNumeric = $true },
    @{ Row = 23; D = '28.014.17'; E = '  +3.26%  '; DNumeric = $false },
    @{ Row = 24; D = '11.24'; E = '  +2.73%  '; DNumeric = $true },
    @{ Row = 25; D = '2.095'; E = '  +0.61%  '; DNumeric = $true },
    @{ Row = 26; D = '156.70'; E = '  +1.65%  '; DNumeric = $true },
    @{ Row = 27; D = '20.61'; E = '  +2.93%  '; DNumeric = $true },
    @{ Row = 28; D = '2.393'; E = '  +14.56%  '; DNumeric = $true },
    @{ Row = 29; D = '2.004.14'; E = '  +3.79%  '; DNumeric = $false },
    @{ Row = 30; D = '122.69'; E = '  +0.85%  '; DNumeric = $true },
    @{ Row = 31; D = '1.148'; E = '  +8.58%  '; DNumeric = $true },
    @{ Row = 32; D = '0.1035'; E = '  +10.74%  '; DNumeric = $true },
    @{ Row = 33; D = '5.765'; E = '  +6.59%  '; DNumeric = $true },
    @{ Row = 34; D = '3.585'; E = '  -1.66%  '; DNumeric = $true },
    @{ Row = 35; D = '0.02303'; E = '  +4.82%  '; DNumeric = $true },
    @{ Row = 36; D = '8.646'; E = '  +14.95%  '; DNumeric = $true },
    @{ Row = 37; D = '0.2111'; E = '  +5.46%  '; DNumeric = $true },
    @{ Row = 38; D = '5.030'; E = '  +4.95%  '; DNumeric = $true },
    @{ Row = 39; D = '0.06055'; E = '  +1.52%  '; DNumeric = $true },
    @{ Row = 40; D = '11.50'; E = '  +4.32%  '; DNumeric = $true },
    @{ Row = 41; D = '0.6293'; E = '  +4.62%  '; DNumeric = $true },
    @{ Row = 42; D = '1.408'; E = '  -1.95%  '; DNumeric = $true },
    @{ Row = 43; D = '0.9975'; E = '  -0.17%  '; DNumeric = $true },
    @{ Row = 44; D = '1.151'; E = '  +4.73%  '; DNumeric = $true },
    @{ Row = 45; D = '13.42'; E = '  +3.85%  '; DNumeric = $true },
    @{ Row = 46; D = '0.5924'; E = '  +4.32%  '; DNumeric = $true },
    @{ Row = 47; D = '3.670'; E = '  +2.31%  '; DNumeric = $true },
    @{ Row = 48; D = '122.49'; E = '  +2.91%  '; DNumeric = $true },
    @{ Row = 49; D = '1.926'; E = '  +3.78%  '; DNumeric = $true },
    @{ Row = 50; D = '1.136'; E = '  +2.72%  '; DNumeric = $true },
    @{ Row = 51; D = '0.06782'; E = '  +1.86%  '; DNumeric = $true }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($r, 4)
        if ($u.DNumeric) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = $u.E
}
